$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new_customer")
$ws2 = $wb.Worksheets.Item("two_contact")

# --- 1. Remove the now-unused "Hyperlink" cell style (mirrors the upstream
#     cleanup that collapses cellStyleXfs/cellStyles back down to just "Normal").
$hyperlinkStyle = $wb.Styles.Item("Hyperlink")
$hyperlinkStyle.Delete()

# --- 2. Populate the role / emailType / type1-3 lookup table (columns G:K).
$ws.Range("G1").Value = "role"
$ws.Range("H1").Value = "emailType"
$ws.Range("I1").Value = "type1"
$ws.Range("J1").Value = "type2"
$ws.Range("K1").Value = "type3"
$ws.Range("G2").Value = "Agent"
$ws.Range("H2").Value = "Business"
$ws.Range("I2").Value = "Home"
$ws.Range("J2").Value = "Cell"
$ws.Range("K2").Value = "Other"
$ws.Range("G3").Value = "CPA"
$ws.Range("H3").Value = "Personal"
$ws.Range("I3").Value = "Fax"
$ws.Range("J3").Value = "Office"
$ws.Range("K3").Value = "Unknown"
$ws.Range("G4").Value = "Clerk"
$ws.Range("H4").Value = "Company"
$ws.Range("I4").Value = "Cell"
$ws.Range("J4").Value = "Fax"
$ws.Range("K4").Value = "Home"
$ws.Range("G5").Value = "Consultant"
$ws.Range("H5").Value = "Government"
$ws.Range("I5").Value = "Office"
$ws.Range("J5").Value = "Cell"
$ws.Range("K5").Value = "Fax"
$ws.Range("G6").Value = "Accountant"
$ws.Range("H6").Value = "Unknown"
$ws.Range("I6").Value = "Other"
$ws.Range("J6").Value = "Unknown"
$ws.Range("K6").Value = "Office"
$ws.Range("G7").Value = "Analyst"
$ws.Range("H7").Value = "Business"
$ws.Range("I7").Value = "Unknown"
$ws.Range("J7").Value = "Home"
$ws.Range("K7").Value = "Cell"

# --- 3. Re-apply the pre-existing cell formats so new cells land on the same
#     style slots as their siblings instead of free-floating new xf records.
$xlPasteFormats = -4122
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial($xlPasteFormats)
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("I1").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("J1").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteFormats)
$ws.Range("G2").Style = "Normal"
$ws.Range("M1").Copy()
$ws.Range("H2").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("J2").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("K2").PasteSpecial($xlPasteFormats)
$ws.Range("G3").Style = "Normal"
$ws.Range("M1").Copy()
$ws.Range("H3").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("J3").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("K3").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("G4").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("H4").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("J4").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("K4").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("G5").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("H5").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("I5").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("K5").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("H6").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("I6").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("J6").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("K6").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("H7").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("I7").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Copy()
$ws.Range("J7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- 4. The handful of cells whose value is one of the brand-new "Cell" /
#     "Accountant" / "Analyst" strings get the distinguishing Consolas style
#     (previously the unused Hyperlink font slot, now repurposed).
$rng = $ws.Range("I4,J5,G6,G7,K7")
$rng.Style = "Normal"
$rng.Font.Name = "Consolas"
$rng.Font.Size = 9
$rng.Font.Color = 2236962

# --- 5. Cursor / selection bookkeeping (cosmetic, matches the saved view state).
$ws.Activate()
$ws.Range("K6").Select()
$ws2.Activate()
$ws2.Range("D20").Select()
$ws.Activate()
